# "ajout surface total  #47"
# Adds a new "Surface totale :" / "Surface" row to the "tableauBatiment"
# table on slide 2, right after the "Adresse" row (i.e. it becomes the
# new 2nd row, pushing "Surface totale chauffee" etc. down by one).

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(2)
$shp = $s.Shapes.Item(1)
$tbl = $shp.Table

# Insert the new row as row 2 (right after "Adresse"); PowerPoint copies
# formatting/height from the following row ("Surface totale chauffee").
$newRow = $tbl.Rows.Add(2)

$tbl.Cell(2, 1).Shape.TextFrame.TextRange.Text = "Surface totale :"
$tbl.Cell(2, 2).Shape.TextFrame.TextRange.Text = "Surface"

# Re-layout the graphic frame so its stored position/extent reflects the
# taller table (one extra row of text).
$shp.Top    = 136.5474
$shp.Height = 199.9467717

# Restore each row's nominal height (the frame re-layout above otherwise
# redistributes the extra space across every row).
$tbl.Rows.Item(1).Height = 15.20535433
$tbl.Rows.Item(2).Height = 39.66015748
$tbl.Rows.Item(3).Height = 39.66015748
$tbl.Rows.Item(4).Height = 31.39007874
$tbl.Rows.Item(5).Height = 31.63637795
